$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- "report" sheet: move selection from C3 to C2 (tabSelected will move off this sheet
#     automatically once the active sheet changes below) ---
$ws1 = $wb.Worksheets.Item("report")
$ws1.Range("C2").Select() | Out-Null

# --- Remove the "MarkIn" sheet entirely (sheet + its shared string + it becomes the
#     deleted reference in WorkDuration's header) ---
$wsMarkIn = $wb.Worksheets.Item("MarkIn")
$wsMarkIn.Delete() | Out-Null

# --- "WorkDuration" sheet: append 8 new "SinglePunchAbsent*" rows (17-24) ---
$ws4 = $wb.Worksheets.Item("WorkDuration")

$newRows = @(
    "SinglePunchAbsent",
    "SinglePunchAbsent Holiday Leave",
    "SinglePunchAbsent Leave Weekoff",
    "SinglePunchAbsent Holiday Weekoff",
    "SinglePunchAbsent Leave ",
    "SinglePunchAbsent Weekoff",
    "SinglePunchAbsent Holiday",
    "SinglePunchAbsent Leave Weekoff Holiday"
)

# Grab the format (Arial 10) already used by column C's "DisplayAll" cells so the new
# rows reuse the same style index instead of registering new ones.
$fmtSrc = $ws4.Cells.Item(2, 3)
$fmtSrc.Copy() | Out-Null

$r = 17
foreach ($val in $newRows) {
    $ws4.Cells.Item($r, 1).Value = $val
    $ws4.Cells.Item($r, 2).Value = $val
    $c3 = $ws4.Cells.Item($r, 3)
    $c3.Value = "DisplayAll"
    $c3.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws4.Cells.Item($r, 4).Value = "Default"
    $ws4.Cells.Item($r, 5).Value = "Default"
    $ws4.Cells.Item($r, 6).Value = $val
    $ws4.Cells.Item($r, 7).Value = "yes"
    $r++
}

# Widen columns A, B and F to fit the longer new text.
$ws4.Columns.Item(1).ColumnWidth = 34.333333333333336
$ws4.Columns.Item(2).ColumnWidth = 36.166666666666664
$ws4.Columns.Item(6).ColumnWidth = 36.0

# Move the selection/active cell to A20 - this also makes WorkDuration the active
# (tabSelected) sheet, matching the post-edit workbook state.
$ws4.Range("A20").Select() | Out-Null
